$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 431
$wsExhibition.Range("F5").Value = 473
$wsExhibition.Range("F7").Value = 2508
$wsExhibition.Range("F9").Value = 6686
$wsExhibition.Range("F10").Value = 177
$wsExhibition.Range("F12").Value = 30

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 431
$wsAll.Range("F5").Value = 473
$wsAll.Range("F9").Value = 2508
$wsAll.Range("F11").Value = 6686
$wsAll.Range("F12").Value = 177
$wsAll.Range("F16").Value = 30
